$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11 (before KETOLAC) for EGYCUSATE, pushing existing
# rows 11-16 (incl. totals/footer) down by one.
$ws.Rows("11:11").Insert()

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "EGYCUSATE 20 MG/5 ML SYRUP 100ML"
$ws.Range("H11").Value = "0:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "25.00"
$ws.Range("P11").Value = "25.0000"
$ws.Range("Q11").Value = "1:0"

# Insert a new row at row 14 (before VOLTAREN, which is now on row 14) for
# URIVIN-N, pushing the remaining rows down by one.
$ws.Rows("14:14").Insert()

$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "URIVIN-N 10 EFF. SACHETS"
$ws.Range("H14").Value = "6:0"
$ws.Range("L14").Value = "1"
$ws.Range("N14").Value = "31.00"
$ws.Range("P14").Value = "31.0000"
$ws.Range("Q14").Value = "0:1"

# Renumber the "index" column for the rows that followed the original
# KETOLAC row (they kept their relative order but shifted down).
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11

# Update the running total shown under the table.
$ws.Range("P18").Value = 515.39

# Update the generated timestamp footer.
$ws.Range("A19").Value = "Wednesday, 23 July, 2025 10:38 AM"
